# Applies the "Corrected excel sheets for application fix issues" edit:
#  - Summary sheet: update Over Due total, Principal/Outstanding for row 3,
#    and touch column G so the used range grows to A1:G5.
#  - Repayment schedule: the (now-irrelevant) first accrued-but-unpaid
#    instalment drops out of the schedule, so every later instalment's
#    day-count/date/figures shift up one row and get recalculated; the
#    now-unused "Over Due" (O) column figures are cleared out too.
#  - Transactions: renumber the two transaction IDs.
#  - Selections/active tab move from "Repayment schedule" to "Summary".

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsTransactions = $wb.Worksheets.Item("Transactions")

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 717.51
$wsSummary.Range("E3").Value = 523.54

# Touch G2 so the sheet's used range extends to column G (matches the new
# <dimension ref="A1:G5"/> / spans="1:7"), leaving it blank afterwards.
$wsSummary.Cells.Item(2, 7).Value = 0
$wsSummary.Cells.Item(2, 7).ClearContents()

# ---------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------

# Row 2: drop the stray "Over Due" figure in column P.
$wsSchedule.Range("P2").Clear()

# Row 3: drop the stray "Over Due" figure in column O.
$wsSchedule.Range("O3").Clear()

# Row 4
$wsSchedule.Range("B4").Value = 31
$wsSchedule.Range("C4").Value = 42095
$wsSchedule.Range("F4").Value = 869.7
$wsSchedule.Range("G4").Value = 8360.5
$wsSchedule.Range("H4").Value = 94.07
$wsSchedule.Range("O4").Clear()

# Row 5
$wsSchedule.Range("B5").Value = 30
$wsSchedule.Range("C5").Value = 42125
$wsSchedule.Range("F5").Value = 881.31
$wsSchedule.Range("G5").Value = 7479.19
$wsSchedule.Range("H5").Value = 82.46
$wsSchedule.Range("O5").Clear()

# Row 6
$wsSchedule.Range("B6").Value = 31
$wsSchedule.Range("C6").Value = 42156
$wsSchedule.Range("F6").Value = 887.54
$wsSchedule.Range("G6").Value = 6591.65
$wsSchedule.Range("H6").Value = 76.23
$wsSchedule.Range("O6").Clear()

# Row 7
$wsSchedule.Range("B7").Value = 30
$wsSchedule.Range("C7").Value = 42186
$wsSchedule.Range("F7").Value = 898.76
$wsSchedule.Range("G7").Value = 5692.89
$wsSchedule.Range("H7").Value = 65.01
$wsSchedule.Range("O7").Clear()

# Row 8
$wsSchedule.Range("B8").Value = 31
$wsSchedule.Range("C8").Value = 42217
$wsSchedule.Range("F8").Value = 905.75
$wsSchedule.Range("G8").NumberFormat = "#,##0.00"
$wsSchedule.Range("G8").Value = 4787.14
$wsSchedule.Range("H8").Value = 58.02
$wsSchedule.Range("O8").Clear()

# Row 9
$wsSchedule.Range("C9").Value = 42248
$wsSchedule.Range("F9").Value = 914.98
$wsSchedule.Range("G9").Value = 3872.16
$wsSchedule.Range("H9").Value = 48.79
$wsSchedule.Range("O9").Clear()

# Row 10
$wsSchedule.Range("B10").Value = 30
$wsSchedule.Range("C10").Value = 42278
$wsSchedule.Range("F10").Value = 925.58
$wsSchedule.Range("G10").Value = 2946.58
$wsSchedule.Range("H10").Value = 38.19
$wsSchedule.Range("O10").Clear()

# Row 11
$wsSchedule.Range("B11").Value = 31
$wsSchedule.Range("C11").Value = 42309
$wsSchedule.Range("F11").Value = 933.74
$wsSchedule.Range("G11").Value = 2012.84
$wsSchedule.Range("H11").Value = 30.03
$wsSchedule.Range("O11").Clear()

# Row 12
$wsSchedule.Range("B12").Value = 30
$wsSchedule.Range("C12").Value = 42339
$wsSchedule.Range("F12").Value = 943.92
$wsSchedule.Range("G12").NumberFormat = "#,##0.00"
$wsSchedule.Range("G12").Value = 1068.92
$wsSchedule.Range("H12").Value = 19.85
$wsSchedule.Range("O12").Clear()

# Row 13 (final instalment)
$wsSchedule.Range("B13").Value = 31
$wsSchedule.Range("C13").Value = 42370
$wsSchedule.Range("F13").NumberFormat = "#,##0.00"
$wsSchedule.Range("F13").Value = 1068.92
$wsSchedule.Range("H13").Value = 10.89
$wsSchedule.Range("K13").NumberFormat = "#,##0.00"
$wsSchedule.Range("K13").Value = 1079.81
$wsSchedule.Range("O13").Clear()
$wsSchedule.Range("P13").NumberFormat = "#,##0.00"
$wsSchedule.Range("P13").Value = 1079.81

# ---------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------
$wsTransactions.Range("A2").Value = 6344
$wsTransactions.Range("A3").Value = 6342

# ---------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------
$wsTransactions.Range("D3").Select()
$wsSchedule.Range("M9").Select()

$wsSummary.Activate()
$wsSummary.Range("C5:D5").Select()
